# Apply the Sep 3 2023 cryptos-list refresh: updated prices/volumes for
# existing rows, plus two coin swaps (rows 37/38 and 50/51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.896.07'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.640.21'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.02'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5046'
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2579'
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06405'
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.68'
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07793'
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = '1.668.45'
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.296'
$ws.Range('E13').Value = '  +1.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5443'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').Value = '0.0₅7873'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.02'
$ws.Range('E16').Value = '  +2.40%  '
$ws.Range('D17').Value = '25.960.00'
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.007'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '198.00'
$ws.Range('E19').Value = '  -2.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.401'
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.989'
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.008'
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.007'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.880'
$ws.Range('E24').Value = '  -3.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '140.26'
$ws.Range('E25').Value = '  -1.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1143'
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.878'
$ws.Range('E27').Value = '  +2.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.75'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05022'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.258'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.201'
$ws.Range('E32').Value = '  +0.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.529'
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.362'
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.8957'
$ws.Range('E35').Value = '  +0.31%  '
$ws.Range('E36').Value = '  -1.31%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.138.60'
$ws.Range('E37').Value = '  -2.47%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5531'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01557'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.007'
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.712'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8163'
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.49'
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').Value = '0.0₈122'
$ws.Range('E44').Value = '  +13.44%  '
$ws.Range('D45').Value = '1.778.13'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4544'
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.35'
$ws.Range('E47').Value = '  +0.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.005'
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05083'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.004'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.09525'
$ws.Range('E51').Value = '  +2.88%  '
